$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 41
$ws.Range("H41").Value = 326.7
$ws.Range("I41").Value = 55
$ws.Range("J41").Value = 394.625
$ws.Range("K41").Value = 55
$ws.Range("L41").Value = 394.625
$ws.Range("M41").Value = 385
$ws.Range("N41").Value = -1274.625

# ALC row 64
$ws.Range("H64").Value = 3306.077
$ws.Range("I64").Value = 3129.8333
$ws.Range("J64").Value = 3457.1428
$ws.Range("K64").Value = 3129.8333
$ws.Range("L64").Value = 3457.1428
$ws.Range("M64").Value = -2881.8333
$ws.Range("N64").Value = -3953.1428

# ALC row 67
$ws.Range("H67").Value = 3306.077
$ws.Range("I67").Value = 3129.8333
$ws.Range("J67").Value = 3457.1428
$ws.Range("K67").Value = 3129.8333
$ws.Range("L67").Value = 3457.1428
$ws.Range("M67").Value = -2271.8333
$ws.Range("N67").Value = -5173.1428

# ALC row 75
$ws.Range("H75").Value = 31876.75
$ws.Range("I75").Value = 14000
$ws.Range("J75").Value = 34430.57
$ws.Range("K75").Value = 14000
$ws.Range("L75").Value = 34430.57
$ws.Range("M75").Value = -13064
$ws.Range("N75").Value = -36302.57

# ALC row 78
$ws.Range("H78").Value = 31876.75
$ws.Range("I78").Value = 14000
$ws.Range("J78").Value = 34430.57
$ws.Range("K78").Value = 42000
$ws.Range("L78").Value = 103291.71
$ws.Range("M78").Value = -37320
$ws.Range("N78").Value = -112651.71

# ALC row 100
$ws.Range("H100").Value = 37123070
$ws.Range("I100").Value = 55558108
$ws.Range("J100").Value = 253003
$ws.Range("K100").Value = 55558108
$ws.Range("L100").Value = 253003
$ws.Range("M100").Value = -55557567
$ws.Range("N100").Value = -254085

# ALC row 106
$ws.Range("H106").Value = 5235
$ws.Range("I106").Value = 6425
$ws.Range("J106").Value = 4937.5
$ws.Range("K106").Value = 6425
$ws.Range("L106").Value = 4937.5
$ws.Range("M106").Value = -5794
$ws.Range("N106").Value = -6199.5

# ALC row 132
$ws.Range("H132").Value = 3526479.2
$ws.Range("I132").Value = 4619385.5
$ws.Range("J132").Value = 4892.5557
$ws.Range("K132").Value = 13858156.5
$ws.Range("L132").Value = 14677.6671
$ws.Range("M132").Value = -13855626.5
$ws.Range("N132").Value = -19737.6671

# ALC row 137
$ws.Range("H137").Value = 1219.0513
$ws.Range("I137").Value = 744.2857
$ws.Range("J137").Value = 2427.5454
$ws.Range("K137").Value = 2232.8571
$ws.Range("L137").Value = 7282.6362
$ws.Range("M137").Value = 317.1428999999998
$ws.Range("N137").Value = -12382.6362

# ALC row 138
$ws.Range("H138").Value = 4607.013
$ws.Range("I138").Value = 1756.1
$ws.Range("J138").Value = 5625.1963
$ws.Range("K138").Value = 5268.299999999999
$ws.Range("L138").Value = 16875.5889
$ws.Range("M138").Value = -128.2999999999993
$ws.Range("N138").Value = -27155.5889

$ws = $wb.Worksheets.Item("ARM")
# ARM row 98
$ws.Range("H98").Value = 34982.5
$ws.Range("J98").Value = 34982.5
$ws.Range("L98").Value = 34982.5
$ws.Range("N98").Value = -40972.5

# ARM row 133
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("BSM")
# BSM row 80
$ws.Range("H80").Value = 805.5
$ws.Range("I80").Value = 255
$ws.Range("J80").Value = 884.1429000000001
$ws.Range("K80").Value = 255
$ws.Range("L80").Value = 884.1429000000001
$ws.Range("M80").Value = 743
$ws.Range("N80").Value = -2880.1429

# BSM row 83
$ws.Range("H83").Value = 805.5
$ws.Range("I83").Value = 255
$ws.Range("J83").Value = 884.1429000000001
$ws.Range("K83").Value = 1275
$ws.Range("L83").Value = 4420.7145
$ws.Range("M83").Value = 3717
$ws.Range("N83").Value = -14404.7145

# BSM row 86
$ws.Range("H86").Value = 2444.6365
$ws.Range("I86").Value = 1985.5
$ws.Range("J86").Value = 3248.125
$ws.Range("K86").Value = 1985.5
$ws.Range("L86").Value = 3248.125
$ws.Range("M86").Value = -862.5
$ws.Range("N86").Value = -5494.125

# BSM row 89
$ws.Range("H89").Value = 2444.6365
$ws.Range("I89").Value = 1985.5
$ws.Range("J89").Value = 3248.125
$ws.Range("K89").Value = 9927.5
$ws.Range("L89").Value = 16240.625
$ws.Range("M89").Value = -4311.5
$ws.Range("N89").Value = -27472.625

# BSM row 107
$ws.Range("H107").Value = 890.8182
$ws.Range("I107").Value = 884.37933
$ws.Range("J107").Value = 937.5
$ws.Range("K107").Value = 884.37933
$ws.Range("L107").Value = 937.5
$ws.Range("M107").Value = 1035.62067
$ws.Range("N107").Value = -4777.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 34570.594
$ws.Range("I16").Value = 51366.7
$ws.Range("J16").Value = 6577.0835
$ws.Range("K16").Value = 51366.7
$ws.Range("L16").Value = 6577.0835
$ws.Range("M16").Value = -51079.7
$ws.Range("N16").Value = -7151.0835

# CRP row 31
$ws.Range("H31").Value = 12503741
$ws.Range("I31").Value = 20834828
$ws.Range("J31").Value = 7111.5
$ws.Range("K31").Value = 20834828
$ws.Range("L31").Value = 7111.5
$ws.Range("M31").Value = -20834533
$ws.Range("N31").Value = -7701.5

# CRP row 34
$ws.Range("H34").Value = 12503741
$ws.Range("I34").Value = 20834828
$ws.Range("J34").Value = 7111.5
$ws.Range("K34").Value = 20834828
$ws.Range("L34").Value = 7111.5
$ws.Range("M34").Value = -20834626
$ws.Range("N34").Value = -7515.5

# CRP row 113
$ws.Range("H113").Value = 34570.594
$ws.Range("I113").Value = 51366.7
$ws.Range("J113").Value = 6577.0835
$ws.Range("K113").Value = 51366.7
$ws.Range("L113").Value = 6577.0835
$ws.Range("M113").Value = -49196.7
$ws.Range("N113").Value = -10917.0835

# CRP row 132
$ws.Range("H132").Value = 2117.3333
$ws.Range("I132").Value = 1498.5294
$ws.Range("J132").Value = 4747.25
$ws.Range("K132").Value = 4495.5882
$ws.Range("L132").Value = 14241.75
$ws.Range("M132").Value = -1965.5882
$ws.Range("N132").Value = -19301.75

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 609.7368
$ws.Range("I5").Value = 456.07144
$ws.Range("J5").Value = 1040
$ws.Range("K5").Value = 1368.21432
$ws.Range("L5").Value = 3120
$ws.Range("M5").Value = -1256.21432
$ws.Range("N5").Value = -3344

# CUL row 64
$ws.Range("H64").Value = 5566.6665
$ws.Range("I64").Value = 5766.6665
$ws.Range("J64").Value = 5466.6665
$ws.Range("K64").Value = 17299.9995
$ws.Range("L64").Value = 16399.9995
$ws.Range("M64").Value = -17029.9995
$ws.Range("N64").Value = -16939.9995

# CUL row 67
$ws.Range("H67").Value = 5566.6665
$ws.Range("I67").Value = 5766.6665
$ws.Range("J67").Value = 5466.6665
$ws.Range("K67").Value = 17299.9995
$ws.Range("L67").Value = 16399.9995
$ws.Range("M67").Value = -16363.9995
$ws.Range("N67").Value = -18271.9995

# CUL row 108
$ws.Range("H108").Value = 725.2
$ws.Range("I108").Value = 725.2
$ws.Range("K108").Value = 2175.6
$ws.Range("M108").Value = 704.3999999999996

# CUL row 113
$ws.Range("H113").Value = 529.6667
$ws.Range("I113").Value = 499
$ws.Range("J113").Value = 545
$ws.Range("K113").Value = 1497
$ws.Range("L113").Value = 1635
$ws.Range("M113").Value = 673
$ws.Range("N113").Value = -5975

# CUL row 121
$ws.Range("H121").Value = 50006584
$ws.Range("I121").Value = 1400
$ws.Range("J121").Value = 55562716
$ws.Range("K121").Value = 4200
$ws.Range("L121").Value = 166688148
$ws.Range("M121").Value = -2890
$ws.Range("N121").Value = -166690768

# CUL row 126
$ws.Range("H126").Value = 2995
$ws.Range("I126").Value = 1997.5
$ws.Range("K126").Value = 5992.5
$ws.Range("M126").Value = -1052.5

# CUL row 131
$ws.Range("H131").Value = 88107
$ws.Range("I131").Value = 328.42856
$ws.Range("J131").Value = 126510.125
$ws.Range("K131").Value = 985.28568
$ws.Range("L131").Value = 379530.375
$ws.Range("M131").Value = 4054.71432
$ws.Range("N131").Value = -389610.375

# CUL row 135
$ws.Range("H135").Value = 609.7368
$ws.Range("I135").Value = 456.07144
$ws.Range("J135").Value = 1040
$ws.Range("K135").Value = 4104.64296
$ws.Range("L135").Value = 9360
$ws.Range("M135").Value = -1569.64296
$ws.Range("N135").Value = -14430

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97
$ws.Range("H97").Value = 977.3889
$ws.Range("I97").Value = 602.5862
$ws.Range("J97").Value = 2530.1428
$ws.Range("K97").Value = 602.5862
$ws.Range("L97").Value = 2530.1428
$ws.Range("M97").Value = -106.5862
$ws.Range("N97").Value = -3522.1428

# GSM row 113
$ws.Range("H113").Value = 6668838.5
$ws.Range("I113").Value = 14287994
$ws.Range("J113").Value = 2077.375
$ws.Range("K113").Value = 14287994
$ws.Range("L113").Value = 2077.375
$ws.Range("M113").Value = -14285824
$ws.Range("N113").Value = -6417.375

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 1431
$ws.Range("I16").Value = 1431
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1431
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1261
$ws.Range("N16").ClearContents()

# LTW row 100
$ws.Range("H100").Value = 2190.5625
$ws.Range("I100").Value = 2096.077
$ws.Range("J100").Value = 2600
$ws.Range("K100").Value = 2096.077
$ws.Range("L100").Value = 2600
$ws.Range("M100").Value = -1555.077
$ws.Range("N100").Value = -3682

$ws = $wb.Worksheets.Item("WVR")
# WVR row 107
$ws.Range("H107").Value = 20838214
$ws.Range("I107").Value = 31250746
$ws.Range("J107").Value = 13147.875
$ws.Range("K107").Value = 93752238
$ws.Range("L107").Value = 39443.625
$ws.Range("M107").Value = -93750318
$ws.Range("N107").Value = -43283.625

# WVR row 113
$ws.Range("H113").Value = 566
$ws.Range("I113").Value = 315.85
$ws.Range("J113").Value = 1066.3
$ws.Range("K113").Value = 947.5500000000001
$ws.Range("L113").Value = 3198.9
$ws.Range("M113").Value = 1222.45
$ws.Range("N113").Value = -7538.9

# WVR row 122
$ws.Range("H122").Value = 1422.2858
$ws.Range("I122").Value = 921.2
$ws.Range("J122").Value = 2675
$ws.Range("K122").Value = 2763.6
$ws.Range("L122").Value = 8025
$ws.Range("M122").Value = -313.6000000000004
$ws.Range("N122").Value = -12925
